$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sprint 1 backlog")
$ws.Range("T10").Value = 0
$ws.Range("U10").Value = 0
$ws.Range("V10").Value = 0
$ws.Range("W10").Value = 0
$ws.Range("T10:W10").Font.Bold = $true

$ws.Range("B10").Font.Bold = $true

# Reflect the new selection/focus on each sheet (cursor ends up on the
# newly-entered cell of the active "Sprint 1 backlog" sheet).
$wsProduct = $wb.Worksheets.Item("Product backlog")
$wsProduct.Activate()
$wsProduct.Range("D4").Select()

$wsSprintTasks = $wb.Worksheets.Item("Sprint tasks")
$wsSprintTasks.Activate()
$wsSprintTasks.Range("E13").Select()

$ws.Activate()
$ws.Range("W10").Select()
